# weiter im TODO. Links repariert
#
# Re-order the slide deck so that the two small "NOT AUS" / "HAUPT SCHALTER"
# slides (originally slides 2 and 3) come right after the "Steuerungspanel"
# slide (originally slide 5), instead of right after the first
# Energie/Stoff/Signal overview slide.
#
# Before: [Energie/Stoff/Signal] [NOT AUS] [HAUPT SCHALTER] [Schaltnetzteil] [Steuerungspanel] [Energie/Stoff/Signal]
# After : [Energie/Stoff/Signal] [Schaltnetzteil] [Steuerungspanel] [NOT AUS] [HAUPT SCHALTER] [Energie/Stoff/Signal]

$p = $ppt.ActivePresentation

# Move the "Schaltnetzteil" slide (slide 4) up to position 2 -- this pushes
# "NOT AUS" and "HAUPT SCHALTER" down by one each.
$p.Slides.Item(4).MoveTo(2)

# The "Steuerungspanel" slide is now at position 5 (it didn't move when the
# slide above moved up); bring it up to position 3, right after
# "Schaltnetzteil".
$p.Slides.Item(5).MoveTo(3)

Write-Output "Final slide order:"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $title = ""
    if ($slide.Shapes.Count -ge 1) {
        try { $title = $slide.Shapes.Item(1).TextFrame.TextRange.Text } catch {}
    }
    Write-Output "  $i : $title"
}
